$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-18 Thursday" "2025-12-19 Friday"

Replace-Text "588×8=4704" "840×8=6720"
Replace-Text "836×6=5016" "630×9=5670"
Replace-Text "499×8=3992" "605×7=4235"
Replace-Text "151×3=453" "953×7=6671"
Replace-Text "590×5=2950" "386×9=3474"
Replace-Text "401×9=3609" "939×4=3756"
Replace-Text "124×8=992" "185×9=1665"
Replace-Text "236×5=1180" "542×5=2710"
Replace-Text "311×2=622" "792×7=5544"
Replace-Text "584×7=4088" "392×7=2744"
Replace-Text "486×3=1458" "675×9=6075"
Replace-Text "479×6=2874" "748×5=3740"
Replace-Text "359×5=1795" "756×3=2268"
Replace-Text "607×6=3642" "115×2=230"
Replace-Text "817×2=1634" "329×4=1316"
Replace-Text "201×6=1206" "471×8=3768"
Replace-Text "248×6=1488" "949×5=4745"
Replace-Text "802×9=7218" "147×3=441"
Replace-Text "794×3=2382" "941×9=8469"
Replace-Text "143×6=858" "427×8=3416"
Replace-Text "252×2=504" "695×4=2780"
Replace-Text "823×7=5761" "880×6=5280"
Replace-Text "510×5=2550" "758×7=5306"
Replace-Text "733×2=1466" "151×8=1208"
Replace-Text "393×6=2358" "642×7=4494"

Write-Host "Done"
